$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "moved against" columns (X3/Y3) for the existing last row (row 3)
$ws.Range("X3").Value = 1.0200049999999976
$ws.Range("Y3").Value = "Up"

# Create new row 4 by copying the formatting of row 3, then overwrite with new data
$ws.Range("A3:W3").Copy($ws.Range("A4:W4"))

$ws.Range("A4").Value = 42641.892569444448
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 14
$ws.Range("E4").Value = 14760
$ws.Range("F4").Value = 2208
$ws.Range("G4").Value = 59
$ws.Range("H4").Value = 37
$ws.Range("I4").Value = 77
$ws.Range("J4").Value = 22
$ws.Range("K4").Value = 15137
$ws.Range("L4").Value = 333
$ws.Range("M4").Value = 211
$ws.Range("N4").Value = 14
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = "Noun"
$ws.Range("Q4").Value = 52.89259217263573
$ws.Range("R4").Value = 0.85
$ws.Range("S4").Value = 0.020199999999999999
$ws.Range("T4").Value = -0.0172
$ws.Range("U4").Value = 15
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0
